{"js": "// Word & grammar corrections (matches the six wording fixes in the diff):\n//   1. \"well known\"  -> \"well-known\"      (hyphenated compound adjective)\n//   2. \"ask\"         -> \"be asked\"        (passive voice correction)\n//   3. \"received from the bank's branch\"  -> \"received from all the bank's branch\"\n//   4. \"base\"        -> \"based\"           (verb agreement)\n//   5. \"First, more robust\" -> \"First, it is more robust\" (missing subject)\n//   6. \"have\"        -> \"has\"             (subject/verb agreement)\n\nconst body = context.document.body;\n\n// Each entry anchors on a short, unique phrase around the words being\n// corrected, then rewrites just that anchor span (via Range.insertText\n// with Word.InsertLocation.replace) so the rest of each paragraph's text\n// and formatting is left untouched.\nconst edits = [\n  {\n    find: \"of well known or being rated\",\n    replace: \"of well-known or being rated\",\n  },\n  {\n    find: \"they will ask how much in terms\",\n    replace: \"they will be asked how much in terms\",\n  },\n  {\n    find: \"the form received from the bank\\u2019s branch)\",\n    replace: \"the form received from all the bank\\u2019s branch)\",\n  },\n  {\n    find: \"It is base on the below formula\",\n    replace: \"It is based on the below formula\",\n  },\n  {\n    find: \"has several advantages. First, more robust to the fluctuation\",\n    replace: \"has several advantages. First, it is more robust to the fluctuation\",\n  },\n  {\n    find: \"Third everybody have the chances to apply the stocks\",\n    replace: \"Third everybody has the chances to apply the stocks\",\n  },\n];\n\nfor (const edit of edits) {\n  const results = body.search(edit.find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${edit.find}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(edit.replace, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word & grammar corrections (matches the six wording fixes in the diff):\n#   1. \"well known\"  -> \"well-known\"      (hyphenated compound adjective)\n#   2. \"ask\"         -> \"be asked\"        (passive voice correction)\n#   3. \"received from the bank's branch\"  -> \"received from all the bank's branch\"\n#   4. \"base\"        -> \"based\"           (verb agreement)\n#   5. \"First, more robust\" -> \"First, it is more robust\" (missing subject)\n#   6. \"have\"        -> \"has\"             (subject/verb agreement)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ Find = \"of well known or being rated\"; Replace = \"of well-known or being rated\" },\n    @{ Find = \"they will ask how much in terms\"; Replace = \"they will be asked how much in terms\" },\n    @{ Find = \"the form received from the bank\" + [char]0x2019 + \"s branch)\"; Replace = \"the form received from all the bank\" + [char]0x2019 + \"s branch)\" },\n    @{ Find = \"It is base on the below formula\"; Replace = \"It is based on the below formula\" },\n    @{ Find = \"has several advantages. First, more robust to the fluctuation\"; Replace = \"has several advantages. First, it is more robust to the fluctuation\" },\n    @{ Find = \"Third everybody have the chances to apply the stocks\"; Replace = \"Third everybody has the chances to apply the stocks\" }\n)\n\nforeach ($edit in $edits) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($edit.Find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $edit.Replace, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Could not find text: $($edit.Find)\"\n    }\n}\n"}
